$wb = $excel.ActiveWorkbook

# Rename "Joslin estimates" to "Henry County estimates"
$wb.Worksheets.Item("Joslin estimates").Name = "Henry County estimates"

# Delete the "Summer 2024 estimates" sheet (no longer used in the manuscript)
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Summer 2024 estimates").Delete()
